$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 = " " (space), B1 = "X", C1 = "Y", D1 = "Angle"
$ws.Cells.Item(1, 1).Value = " "
$ws.Cells.Item(1, 2).Value = "X"
$ws.Cells.Item(1, 3).Value = "Y"
$ws.Cells.Item(1, 4).Value = "Angle"

# Data rows 2..36 (35 minutiae entries)
$data = @(
    @(1, 71.429000000000002, 47.347000000000001, 201.251),
    @(2, 115.91800000000001, 60, 352.23500000000001),
    @(3, 182.857, 61.633000000000003, 524.93200000000002),
    @(4, 49.387999999999998, 94.694000000000003, 209.05500000000001),
    @(5, 135.51, 93.468999999999994, 177.614),
    @(6, 93.061000000000007, 102.857, 11.31),
    @(7, 55.51, 119.184, 210.06899999999999),
    @(8, 69.796000000000006, 143.673, 9.4619999999999997),
    @(9, 38.776000000000003, 164.898, 225),
    @(10, 84.897999999999996, 165.30600000000001, 204.22800000000001),
    @(11, 195.91800000000001, 164.08199999999999, 324.16199999999998),
    @(12, 161.63300000000001, 172.245, 341.565),
    @(13, 181.63300000000001, 180.816, 159.44399999999999),
    @(14, 57.143000000000001, 189.79599999999999, 40.914000000000001),
    @(15, 167.755, 203.26499999999999, 334.654),
    @(16, 220.816, 233.06100000000001, 147.995),
    @(17, 171.429, 223.673, 158.749),
    @(18, 144.898, 230.61199999999999, 341.565),
    @(19, 89.796000000000006, 222.041, 208.072),
    @(20, 69.796000000000006, 233.46899999999999, 37.569000000000003),
    @(21, 46.122, 258.77600000000001, 243.435),
    @(22, 102.041, 240.816, 21.800999999999998),
    @(23, 140, 237.143, 180),
    @(24, 178.77600000000001, 248.16300000000001, 323.13),
    @(25, 197.959, 269.79599999999999, 308.65999999999997),
    @(26, 132.245, 251.429, 340.71),
    @(27, 155.91800000000001, 271.02, 145.49100000000001),
    @(28, 47.755000000000003, 311.02, 56.31),
    @(29, 87.754999999999995, 283.673, 35.537999999999997),
    @(30, 164.08199999999999, 325.714, 129.47200000000001),
    @(31, 107.755, 300.81599999999997, 348.11099999999999),
    @(32, 43.265000000000001, 350.61200000000002, 56.31),
    @(33, 247.755, 272.65300000000002, 138.81399999999999),
    @(34, 209.38800000000001, 371.02, 154.44),
    @(35, 193.87799999999999, 384.49, 146.88900000000001)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the selection to match the post-edit state
[void]$ws.Range("H37").Select()
